$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The schedule's final two rows used to be the "Exam" session (row 43/44 in
# the numbering column, Excel rows 44/45). The new schedule pushes the exam
# back by three extra "Mini Project" sessions, so:
#   - rows 44 and 45 become "Mini Project" sessions (like the rows above them)
#   - three brand new rows (46, 47, 48) are appended
#   - the last of those new rows (48) carries the "Exam" session that used
#     to live in rows 44/45.

# Row 44: Exam -> Mini Project (column A and E stay the same).
$ws.Cells.Item(44, 2).Value = "Mini Project"
$ws.Cells.Item(44, 3).Value = "Practical"
$ws.Cells.Item(44, 4).Value = "Project work"

# Row 45: Exam -> Mini Project (column A and E stay the same).
$ws.Cells.Item(45, 2).Value = "Mini Project"
$ws.Cells.Item(45, 3).Value = "Practical"
$ws.Cells.Item(45, 4).Value = "Project work"

# Append three new rows at the bottom of the schedule.
$ws.Rows.Item(46).Insert()
$ws.Rows.Item(47).Insert()
$ws.Rows.Item(48).Insert()

# Row 46: another "Mini Project" session.
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "Mini Project"
$ws.Cells.Item(46, 3).Value = "Practical"
$ws.Cells.Item(46, 4).Value = "Project work"
$ws.Cells.Item(46, 5).Value = "OJ"

# Row 47: another "Mini Project" session.
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = "Mini Project"
$ws.Cells.Item(47, 3).Value = "Practical"
$ws.Cells.Item(47, 4).Value = "Project work"
$ws.Cells.Item(47, 5).Value = "OJ"

# Row 48: the "Exam" session, now moved to the very end of the schedule.
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = "Exam"
$ws.Cells.Item(48, 3).Value = "Exam"
$ws.Cells.Item(48, 4).Value = "MCQ Exam"
$ws.Cells.Item(48, 5).Value = "OJ"

# Scroll the view down and select B46, matching where the author was
# working when the new rows were added.
$excel.ActiveWindow.ScrollRow = 33
$ws.Range("B46").Select()
